$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-RowValues($range, $values) {
    $arr = New-Object 'object[,]' 1,$values.Length
    for ($i = 0; $i -lt $values.Length; $i++) {
        $arr[0,$i] = $values[$i]
    }
    $ws.Range($range).Value = $arr
}

Set-RowValues "T2:Y2" @(1, 9, 18, 30, 34, 40)
Set-RowValues "T3:Y3" @(2, 10, 15, 27, 35, 38)
Set-RowValues "T4:Y4" @(7, 13, 28, 29, 30, 38)
Set-RowValues "T5:Y5" @(2, 11, 19, 20, 38, 40)
Set-RowValues "T6:Y6" @(6, 15, 20, 28, 32, 43)
Set-RowValues "K7:Q7" @(2, 13, 25, 28, 29, 36, 34)

$ws.Range("O20").Select()
